$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Salinity_Q_river_splitted"
